# Replace the static "www.drpaulduenas.com" footer text with a
# MERGEFIELD ("=website") construct, matching the other dynamic
# fields already used throughout the footer (address, phone, city, ...).
#
# The target run lives in the *default* (primary) footer story, not in
# Document.Content (which is body-only), so we have to reach it via
# Sections/Footers and operate on the Range object that Find.Execute
# narrows in place.

$d = $word.ActiveDocument

$targetText = "www.drpaulduenas.com"
$fieldName = "website"

$found = $false

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $section = $d.Sections.Item($s)
    $footers = $section.Footers
    for ($f = 1; $f -le $footers.Count; $f++) {
        $ftr = $footers.Item($f)
        if (-not $ftr.Exists) {
            continue
        }

        $rng = $ftr.Range
        $rng.Find.ClearFormatting()
        $rng.Find.Execute($targetText, $true, $false, $false, $false, $false,
                           $true, 1, $false, "", 0)

        if ($rng.Find.Found) {
            $found = $true

            # Clear the matched text in place (this collapses the range
            # to a single insertion point right where the text was).
            $rng.Text = ""

            # Re-insert the same run properties, but as a MERGEFIELD
            # field-code construct (fldChar begin/separate/end +
            # instrText + the cached display text) instead of literal
            # text, mirroring the other merge fields in this footer.
            $instr = " MERGEFIELD =$fieldName \* MERGEFORMAT "
            $display = [char]0x00AB + "=$fieldName" + [char]0x00BB

            $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
                   '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
                   '<pkg:xmlData>' +
                   '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
                   '<w:body>' +
                   '<w:p w14:paraId="24EA949D" w14:textId="77777777" w:rsidR="004D2A29" w:rsidRDefault="004D2A29" w:rsidP="004D2A29">' +
                   '<w:pPr><w:pStyle w:val="Footer"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr>' +
                   '<w:r><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r>' +
                   '<w:r><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:instrText xml:space="preserve">' + $instr + '</w:instrText></w:r>' +
                   '<w:r><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r>' +
                   '<w:r><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>' + $display + '</w:t></w:r>' +
                   '<w:r><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r>' +
                   '</w:p>' +
                   '</w:body></w:document>' +
                   '</pkg:xmlData></pkg:part></pkg:package>'

            $rng.InsertXML($xml)
        }
    }
}

if (-not $found) {
    throw "Could not find '$targetText' in any footer to convert into a MERGEFIELD"
}

Write-Output "Converted '$targetText' into a =$fieldName MERGEFIELD"
